# Applies the "Command as a Concept" planning-doc trivial edits:
#   1. Title paragraph: drop the comma after "Circle Language Spec Plan"
#      and after "... Spec"; split "2008-06 " into two runs and move the
#      "_GoBack" bookmark here (it used to sit in front of the "Goal"
#      heading).
#   2. The four "date" smart tags get their <w:attr> children reordered
#      to Month/Day/Year.
#   3. The "Goal" heading no longer carries the (now relocated) bookmark.
#   4. Heading 2 / Heading 3 styles switch from Arial to Calibri (Heading 2
#      also grows from 16pt to 18pt).
#
# Whole paragraphs are rewritten via Range.InsertXML (wrapped in the
# mandatory pkg:package/pkg:part envelope) so run-splits, the bookmark
# placement and the smart-tag attribute order come out exactly as
# intended, then the two heading styles are restyled through the normal
# Style/Font object model.

$d = $word.ActiveDocument

# NOTE: this interpreter only honours *positional* arguments reliably,
# so this helper is always invoked as `Set-ParagraphXml $para $xml`.
function Set-ParagraphXml {
    param(
        [object]$Paragraph,
        [string]$InnerXml
    )
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData>' `
        + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
        + '<w:body>' + $InnerXml + '</w:body>' `
        + '</w:document>' `
        + '</pkg:xmlData></pkg:part></pkg:package>'
    $Paragraph.Range.InsertXML($xmlFrag)
}

# Locate the paragraphs we need to touch by their (still unique) text so
# the script does not depend on brittle paragraph indices.
$titlePara = $null
$dateRangePara = $null
$goalPara = $null
$timeRangePara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Circle Language Spec Plan*") { $titlePara = $p }
    elseif ($t -like "Date:*") { $dateRangePara = $p }
    elseif ($t -eq "Goal`r") { $goalPara = $p }
    elseif ($t -like "3 months and 1 week*") { $timeRangePara = $p.Previous() }
}

# 1. Title heading paragraph.
$titleXml = '<w:p w:rsidR="00AB55BB" w:rsidRDefault="0071318D" w:rsidP="00AB55BB">' `
    + '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' `
    + '<w:r w:rsidRPr="0071318D"><w:t>Circle Language Spec Plan</w:t></w:r>' `
    + '<w:r w:rsidR="00CF03FF"><w:br/></w:r>' `
    + '<w:r w:rsidR="00D73F50" w:rsidRPr="00D73F50"><w:t>2008-06</w:t></w:r>' `
    + '<w:r w:rsidR="00D73F50" w:rsidRPr="00D73F50"><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' `
    + '<w:r w:rsidR="000D32B0"><w:t>Command</w:t></w:r>' `
    + '<w:r w:rsidR="0081141D"><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t>a</w:t></w:r>' `
    + '<w:r w:rsidR="00290BC1"><w:t xml:space="preserve">s </w:t></w:r>' `
    + '<w:r><w:t>a</w:t></w:r>' `
    + '<w:r w:rsidR="00290BC1"><w:t xml:space="preserve"> Concept</w:t></w:r>' `
    + '<w:r w:rsidR="00313AFB"><w:t xml:space="preserve"> Spec</w:t></w:r>' `
    + '<w:r w:rsidR="00AB55BB"><w:br/><w:t>Project Summary</w:t></w:r>' `
    + '</w:p>'
Set-ParagraphXml $titlePara $titleXml

# 2. "Date: <smartTag> - <smartTag>" paragraph - reorder Month/Day/Year.
$dateXml = '<w:p w:rsidR="00CF03FF" w:rsidRPr="00E776B2" w:rsidRDefault="00CF03FF" w:rsidP="00FF5F4B">' `
    + '<w:pPr><w:ind w:left="284"/><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr></w:pPr>' `
    + '<w:r w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve">Date: </w:t></w:r>' `
    + '<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">' `
    + '<w:smartTagPr><w:attr w:name="Month" w:val="6"/><w:attr w:name="Day" w:val="21"/><w:attr w:name="Year" w:val="2008"/></w:smartTagPr>' `
    + '<w:r w:rsidR="00C26530" w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>June 21, 2008</w:t></w:r>' `
    + '</w:smartTag>' `
    + '<w:r w:rsidR="00C26530" w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve"> &#8211; </w:t></w:r>' `
    + '<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">' `
    + '<w:smartTagPr><w:attr w:name="Month" w:val="12"/><w:attr w:name="Day" w:val="23"/><w:attr w:name="Year" w:val="2008"/></w:smartTagPr>' `
    + '<w:r w:rsidR="000D32B0"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>December 23</w:t></w:r>' `
    + '<w:r w:rsidR="00C26530" w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>, 2008</w:t></w:r>' `
    + '</w:smartTag>' `
    + '</w:p>'
Set-ParagraphXml $dateRangePara $dateXml

# 3. "Goal" heading - the _GoBack bookmark moved up into the title, so it
#    no longer belongs here.
$goalXml = '<w:p w:rsidR="00EF400A" w:rsidRDefault="00EF400A" w:rsidP="00FF5F4B">' `
    + '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' `
    + '<w:r><w:t>Goa</w:t></w:r>' `
    + '<w:r w:rsidR="00FF5F4B"><w:t>l</w:t></w:r>' `
    + '</w:p>'
Set-ParagraphXml $goalPara $goalXml

# 4. "June 23, 2008 - October 1, 2008" paragraph - reorder Month/Day/Year.
$timeXml = '<w:p w:rsidR="00FF5F4B" w:rsidRDefault="00B8645A" w:rsidP="001A67E0">' `
    + '<w:pPr><w:ind w:left="426"/></w:pPr>' `
    + '<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">' `
    + '<w:smartTagPr><w:attr w:name="Month" w:val="6"/><w:attr w:name="Day" w:val="23"/><w:attr w:name="Year" w:val="2008"/></w:smartTagPr>' `
    + '<w:r><w:t>June 23</w:t></w:r>' `
    + '<w:r w:rsidR="00916C2E"><w:t xml:space="preserve">, </w:t></w:r>' `
    + '<w:r w:rsidR="00FF5F4B"><w:t>2008</w:t></w:r>' `
    + '</w:smartTag>' `
    + '<w:r w:rsidR="00916C2E"><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r w:rsidR="00FF5F4B"><w:t xml:space="preserve">&#8211; </w:t></w:r>' `
    + '<w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date">' `
    + '<w:smartTagPr><w:attr w:name="Month" w:val="10"/><w:attr w:name="Day" w:val="1"/><w:attr w:name="Year" w:val="2008"/></w:smartTagPr>' `
    + '<w:r w:rsidR="00A4196C"><w:t>October 1</w:t></w:r>' `
    + '<w:r w:rsidR="00916C2E"><w:t>, 2008</w:t></w:r>' `
    + '</w:smartTag>' `
    + '</w:p>'
Set-ParagraphXml $timeRangePara $timeXml

# 5. Heading 2 / Heading 3 styles: Arial -> Calibri, Heading 2 grows to 18pt.
$heading2 = $d.Styles("Heading 2")
$heading2.Font.Name = "Calibri"
$heading2.Font.Size = 18

$heading3 = $d.Styles("Heading 3")
$heading3.Font.Name = "Calibri"

Write-Output "Applied Command-as-a-Concept trivial edits."
